$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3648
$ws1.Range("F5").Value = 2235
$ws1.Range("F9").Value = 94
$ws1.Range("F11").Value = 1349
$ws1.Range("F13").Value = 2012

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3648
$ws4.Range("F5").Value = 2235
$ws4.Range("F10").Value = 94
$ws4.Range("F13").Value = 3
$ws4.Range("F14").Value = 1349
$ws4.Range("F16").Value = 2012
